$wb = $excel.ActiveWorkbook

# Update "展览" sheet (sheetId=1 / rId1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 8376
$ws1.Range("F5").Value = 354

# Update "全部类型" sheet (sheetId=4 / rId4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 8376
$ws4.Range("F5").Value = 354
